$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.887828350067139
$ws.Range("B1").Value = 2.611997604370117
$ws.Range("C1").Value = 3.527505397796631
$ws.Range("D1").Value = 1.143050312995911
$ws.Range("E1").Value = 0.7340084910392761
